$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 7
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()

# Row 8
$ws.Range("H8").Value = 35
$ws.Range("I8").Value = 31.428572
$ws.Range("J8").Value = 60
$ws.Range("K8").Value = 94.28571599999999
$ws.Range("L8").Value = 180
$ws.Range("M8").Value = 44.71428400000001
$ws.Range("N8").Value = -458

# Row 14
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()

# Row 33
$ws.Range("H33").Value = 16703.2
$ws.Range("I33").Value = 20379.188
$ws.Range("K33").Value = 20379.188
$ws.Range("M33").Value = -20150.188

# Row 53
$ws.Range("H53").Value = 67095.164
$ws.Range("J53").Value = 95493.664
$ws.Range("L53").Value = 95493.664
$ws.Range("N53").Value = -96767.664

# Row 86
$ws.Range("H86").Value = 10966.667
$ws.Range("I86").Value = 9000
$ws.Range("K86").Value = 9000
$ws.Range("M86").Value = -7877

# Row 89
$ws.Range("H89").Value = 10966.667
$ws.Range("I89").Value = 9000
$ws.Range("K89").Value = 45000
$ws.Range("M89").Value = -39384

# Row 101
$ws.Range("H101").Value = 3991.7856
$ws.Range("I101").Value = 820.2
$ws.Range("J101").Value = 5753.778
$ws.Range("K101").Value = 2460.6
$ws.Range("L101").Value = 17261.334
$ws.Range("M101").Value = -838.6000000000004
$ws.Range("N101").Value = -20505.334

# Row 113
$ws.Range("H113").Value = 142859180
$ws.Range("I113").Value = 142859180
$ws.Range("K113").Value = 142859180
$ws.Range("M113").Value = -142855926

# Row 129
$ws.Range("H129").Value = 1264640.8
$ws.Range("I129").Value = 1412
$ws.Range("K129").Value = 4236
$ws.Range("M129").Value = 764

# Row 135
$ws.Range("H135").Value = 2312.6
$ws.Range("I135").Value = 1674.0834
$ws.Range("K135").Value = 15066.7506
$ws.Range("M135").Value = -12531.7506

# Row 137
$ws.Range("H137").Value = 2358.7144
$ws.Range("I137").Value = 2238.2727
$ws.Range("K137").Value = 6714.8181
$ws.Range("M137").Value = -4164.8181

# Row 138
$ws.Range("H138").Value = 3758.1428
$ws.Range("J138").Value = 4521.5
$ws.Range("L138").Value = 13564.5
$ws.Range("N138").Value = -23844.5

# Row 141
$ws.Range("H141").Value = 3140.6086
$ws.Range("I141").Value = 3140.6086
$ws.Range("K141").Value = 9421.825800000001
$ws.Range("M141").Value = -4241.825800000001

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2930.48
$ws.Range("I32").Value = 2785.5857
$ws.Range("K32").Value = 2785.5857
$ws.Range("M32").Value = -2498.5857

# Row 50
$ws.Range("H50").Value = 16396
$ws.Range("J50").Value = 49996
$ws.Range("L50").Value = 49996
$ws.Range("N50").Value = -51424

# Row 61
$ws.Range("H61").Value = 7252996.5
$ws.Range("I61").Value = 8338696
$ws.Range("K61").Value = 8338696
$ws.Range("M61").Value = -8338484

# Row 131
$ws.Range("H131").Value = 69999
$ws.Range("J131").Value = 69999
$ws.Range("L131").Value = 69999
$ws.Range("N131").Value = -80079

# Row 132
$ws.Range("H132").Value = 5540.3335
$ws.Range("I132").Value = 4820.3335
$ws.Range("K132").Value = 14461.0005
$ws.Range("M132").Value = -11931.0005

# Row 136
$ws.Range("H136").Value = 7252996.5
$ws.Range("I136").Value = 8338696
$ws.Range("K136").Value = 25016088
$ws.Range("M136").Value = -25013538

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 4433.75
$ws.Range("I134").Value = 4605.8945
$ws.Range("K134").Value = 13817.6835
$ws.Range("M134").Value = -11282.6835

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5405.143
$ws.Range("I31").Value = 4625.8
$ws.Range("J31").Value = 6304.385
$ws.Range("K31").Value = 4625.8
$ws.Range("L31").Value = 6304.385
$ws.Range("M31").Value = -4330.8
$ws.Range("N31").Value = -6894.385

# Row 34
$ws.Range("H34").Value = 5405.143
$ws.Range("I34").Value = 4625.8
$ws.Range("J34").Value = 6304.385
$ws.Range("K34").Value = 4625.8
$ws.Range("L34").Value = 6304.385
$ws.Range("M34").Value = -4423.8
$ws.Range("N34").Value = -6708.385

# Row 58
$ws.Range("H58").Value = 11899.8
$ws.Range("I58").Value = 7500
$ws.Range("K58").Value = 7500
$ws.Range("M58").Value = -7297

# Row 62
$ws.Range("H62").Value = 5750.067
$ws.Range("I62").Value = 6020.4165
$ws.Range("K62").Value = 6020.4165
$ws.Range("M62").Value = -5396.4165

# Row 65
$ws.Range("H65").Value = 5750.067
$ws.Range("I65").Value = 6020.4165
$ws.Range("K65").Value = 30102.0825
$ws.Range("M65").Value = -26982.0825

# Row 86
$ws.Range("H86").Value = 4822.467
$ws.Range("I86").Value = 3879.7693
$ws.Range("J86").Value = 10950
$ws.Range("K86").Value = 3879.7693
$ws.Range("L86").Value = 10950
$ws.Range("M86").Value = -2756.7693
$ws.Range("N86").Value = -13196

# Row 89
$ws.Range("H89").Value = 4822.467
$ws.Range("I89").Value = 3879.7693
$ws.Range("J89").Value = 10950
$ws.Range("K89").Value = 19398.8465
$ws.Range("L89").Value = 54750
$ws.Range("M89").Value = -13782.8465
$ws.Range("N89").Value = -65982

# Row 105
$ws.Range("H105").Value = 2140.7896
$ws.Range("I105").Value = 2019.1177
$ws.Range("K105").Value = 2019.1177
$ws.Range("M105").Value = -272.1177

# Row 132
$ws.Range("H132").Value = 3742.3333
$ws.Range("I132").Value = 3369.5
$ws.Range("J132").Value = 5047.25
$ws.Range("K132").Value = 10108.5
$ws.Range("L132").Value = 15141.75
$ws.Range("M132").Value = -7578.5
$ws.Range("N132").Value = -20201.75

# Row 134
$ws.Range("H134").Value = 7004.0386
$ws.Range("J134").Value = 9669.467000000001
$ws.Range("L134").Value = 29008.401
$ws.Range("N134").Value = -34078.401

# Row 136
$ws.Range("H136").Value = 11899.8
$ws.Range("I136").Value = 7500
$ws.Range("K136").Value = 22500
$ws.Range("M136").Value = -19950

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 23411384
$ws.Range("I4").Value = 28163538
$ws.Range("J4").Value = 2501899.5
$ws.Range("K4").Value = 84490614
$ws.Range("L4").Value = 7505698.5
$ws.Range("M4").Value = -84490502
$ws.Range("N4").Value = -7505922.5

# Row 5
$ws.Range("H5").Value = 1414.125
$ws.Range("I5").Value = 338.36365
$ws.Range("K5").Value = 1015.09095
$ws.Range("M5").Value = -903.09095

# Row 9
$ws.Range("H9").Value = 4430186
$ws.Range("J9").Value = 467.33334
$ws.Range("L9").Value = 1402.00002
$ws.Range("N9").Value = -1850.00002

# Row 129
$ws.Range("H129").Value = 4248.222
$ws.Range("I129").Value = 2859.875
$ws.Range("K129").Value = 8579.625
$ws.Range("M129").Value = -3579.625

# Row 135
$ws.Range("H135").Value = 1414.125
$ws.Range("I135").Value = 338.36365
$ws.Range("K135").Value = 3045.27285
$ws.Range("M135").Value = -510.2728500000003

# Row 137
$ws.Range("H137").Value = 23772.133
$ws.Range("I137").Value = 2015.3334
$ws.Range("K137").Value = 6046.0002
$ws.Range("M137").Value = -946.0002000000004

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 4583.08
$ws.Range("I2").Value = 235.57895
$ws.Range("J2").Value = 18350.166
$ws.Range("K2").Value = 235.57895
$ws.Range("L2").Value = 18350.166
$ws.Range("M2").Value = -122.57895
$ws.Range("N2").Value = -18576.166

# Row 97
$ws.Range("H97").Value = 837.04
$ws.Range("J97").Value = 687.4286
$ws.Range("L97").Value = 687.4286
$ws.Range("N97").Value = -1679.4286

# Row 99
$ws.Range("H99").Value = 21199.5
$ws.Range("I99").Value = 17439.4
$ws.Range("K99").Value = 17439.4
$ws.Range("M99").Value = -15193.4

# Row 132
$ws.Range("H132").Value = 5130.522
$ws.Range("I132").Value = 4981.048
$ws.Range("K132").Value = 14943.144
$ws.Range("M132").Value = -12413.144

$ws = $wb.Worksheets.Item("LTW")
# Row 133
$ws.Range("H133").Value = 73663.336
$ws.Range("J133").Value = 73663.336
$ws.Range("L133").Value = 73663.336
$ws.Range("N133").Value = -78723.336

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 2814
$ws.Range("I2").Value = 2566
$ws.Range("K2").Value = 2566
$ws.Range("M2").Value = -2454

# Row 4
$ws.Range("H4").Value = 463111.16
$ws.Range("I4").Value = 679166
$ws.Range("J4").Value = 31001.5
$ws.Range("K4").Value = 679166
$ws.Range("L4").Value = 31001.5
$ws.Range("M4").Value = -679053
$ws.Range("N4").Value = -31227.5

# Row 113
$ws.Range("H113").Value = 1111.28
$ws.Range("I113").Value = 1408.5714
$ws.Range("K113").Value = 4225.7142
$ws.Range("M113").Value = -2055.7142

# Row 122
$ws.Range("H122").Value = 5628.4287
$ws.Range("J122").Value = 7933
$ws.Range("L122").Value = 23799
$ws.Range("N122").Value = -28699

# Row 126
$ws.Range("H126").Value = 3801.6667
$ws.Range("I126").Value = 2901.9375
$ws.Range("J126").Value = 10999.5
$ws.Range("K126").Value = 8705.8125
$ws.Range("L126").Value = 32998.5
$ws.Range("M126").Value = -6235.8125
$ws.Range("N126").Value = -37938.5

# Row 135
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
